$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric "Model No" column values ---
$ws.Range("F9").Value = 1
$ws.Range("F11").Value = 2
$ws.Range("F13").Value = 3
$ws.Range("F15").Value = 4
$ws.Range("F17").Value = 5
$ws.Range("F19").Value = 6

# --- Fill text in the order the strings first appear so shared-string table order matches ---
$ws.Range("G9").Value = "MaskRCNN"
$ws.Range("I9").Value = "Not completed "
$ws.Range("H9").Value = " "

$ws.Range("G11").Value = "Base: ResNet`nObjectDet : Retina Net(One Shot)"

$ws.Range("G13").Value = "Baseline Model"
$ws.Range("H13").Value = "Accuracy :0.38`nIoU: .89"
$ws.Range("I13").Value = "Completed but poor scores"

$ws.Range("K7").Value = "Owner"
$ws.Range("K9").Value = "Poornima"
$ws.Range("K11").Value = "Alok"
$ws.Range("K13").Value = "Radhika"

$ws.Range("H15").Value = "IoU:.75`nAccuracy ?:Cannot read"

$ws.Range("K17").Value = "Ganesh"
$ws.Range("G17").Value = "YOLOV4"
$ws.Range("G15").Value = "YOLOV3"
$ws.Range("G19").Value = "SSD"
$ws.Range("K19").Value = "Sachen"

# --- Remaining repeated values (reuse existing shared strings, no new entries) ---
$ws.Range("I11").Value = "Not completed "
$ws.Range("I15").Value = "Completed but poor scores"
$ws.Range("I17").Value = "Not completed "
$ws.Range("K15").Value = "Radhika"

# --- Header row updates (also syncs the Table1 column names) ---
$ws.Range("H7").Value = "Accuracy/IOU Score"

# --- Wrap text + row heights for the multi-line cells ---
$ws.Range("G11").WrapText = $true
$ws.Rows(11).RowHeight = 60

$ws.Range("H13").WrapText = $true
$ws.Rows(13).RowHeight = 30

$ws.Range("H15").WrapText = $true
$ws.Rows(15).RowHeight = 30

# --- Column widths (nearest values this engine's width quantization can reach
#     to the authored widths of 28.5703125 / 24.7109375 / 13.5703125 character units) ---
$ws.Range("I1").ColumnWidth = 27.666666666666664
$ws.Range("J1").ColumnWidth = 23.833333333333336
$ws.Range("K1").ColumnWidth = 12.666666666666668

# --- Selection ---
$ws.Range("H19").Select() | Out-Null
